# Auto-generated Excel COM-interop edit script
# Applies numeric updates to columns H-N across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# matching the profit/price recalculation described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10 (item id G10=1959)
$ws.Range("H10").Value = 19799.334
$ws.Range("J10").Value = 19799.334
$ws.Range("L10").Value = 19799.334
$ws.Range("N10").Value = -20385.334
# Row 33 (item id G33=5512)
$ws.Range("H33").Value = 778.29266
$ws.Range("I33").Value = 604.40625
$ws.Range("J33").Value = 1396.5555
$ws.Range("K33").Value = 604.40625
$ws.Range("L33").Value = 1396.5555
$ws.Range("M33").Value = -375.40625
$ws.Range("N33").Value = -1854.5555
# Row 70 (item id G70=12604)
$ws.Range("H70").Value = 2723.2222
$ws.Range("I70").Value = 5800
$ws.Range("J70").Value = 1844.1428
$ws.Range("K70").Value = 17400
$ws.Range("L70").Value = 5532.428400000001
$ws.Range("M70").Value = -17130
$ws.Range("N70").Value = -6072.428400000001
# Row 73 (item id G73=12604)
$ws.Range("H73").Value = 2723.2222
$ws.Range("I73").Value = 5800
$ws.Range("J73").Value = 1844.1428
$ws.Range("K73").Value = 17400
$ws.Range("L73").Value = 5532.428400000001
$ws.Range("M73").Value = -16464
$ws.Range("N73").Value = -7404.428400000001
# Row 86 (item id G86=12603)
$ws.Range("H86").Value = 9416.875
$ws.Range("I86").Value = 18944.334
$ws.Range("J86").Value = 3700.4
$ws.Range("K86").Value = 18944.334
$ws.Range("L86").Value = 3700.4
$ws.Range("M86").Value = -17821.334
$ws.Range("N86").Value = -5946.4
# Row 89 (item id G89=12603)
$ws.Range("H89").Value = 9416.875
$ws.Range("I89").Value = 18944.334
$ws.Range("J89").Value = 3700.4
$ws.Range("K89").Value = 94721.67
$ws.Range("L89").Value = 18502
$ws.Range("M89").Value = -89105.67
$ws.Range("N89").Value = -29734
# Row 92 (item id G92=19901)
$ws.Range("H92").Value = 9259850
$ws.Range("I92").Value = 18518906
$ws.Range("J92").Value = 793.3333
$ws.Range("K92").Value = 18518906
$ws.Range("L92").Value = 793.3333
$ws.Range("M92").Value = -18517658
$ws.Range("N92").Value = -3289.3333
# Row 106 (item id G106=19903)
$ws.Range("H106").Value = 3499.2307
$ws.Range("I106").Value = 3771.818
$ws.Range("K106").Value = 3771.818
$ws.Range("M106").Value = -3140.818
# Row 112 (item id G112=27960)
$ws.Range("H112").Value = 1213.2667
$ws.Range("J112").Value = 1338.3846
$ws.Range("L112").Value = 4015.1538
$ws.Range("N112").Value = -6231.1538
# Row 129 (item id G129=36115)
$ws.Range("H129").Value = 1098.3077
$ws.Range("J129").Value = 1279.619
$ws.Range("L129").Value = 3838.857
$ws.Range("N129").Value = -13838.857
# Row 138 (item id G138=44169)
$ws.Range("H138").Value = 2199.7673
$ws.Range("I138").Value = 731.52
$ws.Range("J138").Value = 4239
$ws.Range("K138").Value = 2194.56
$ws.Range("L138").Value = 12717
$ws.Range("M138").Value = 2945.44
$ws.Range("N138").Value = -22997

$ws = $wb.Worksheets.Item("ARM")
# Row 25 (item id G25=2471)
$ws.Range("H25").Value = 929
$ws.Range("I25").Value = 905.3333
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 905.3333
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -503.3333
$ws.Range("N25").Value = -1804
# Row 32 (item id G32=44147)
$ws.Range("H32").Value = 8080.44
$ws.Range("I32").Value = 2299.3816
$ws.Range("J32").Value = 26387.125
$ws.Range("K32").Value = 2299.3816
$ws.Range("L32").Value = 26387.125
$ws.Range("M32").Value = -2012.3816
$ws.Range("N32").Value = -26961.125
# Row 110 (item id G110=27708)
$ws.Range("H110").Value = 16632.5
$ws.Range("I110").Value = 20510.166
$ws.Range("K110").Value = 20510.166
$ws.Range("M110").Value = -18465.166

$ws = $wb.Worksheets.Item("BSM")
# Row 10 (item id G10=2417)
$ws.Range("H10").Value = 1200
$ws.Range("J10").Value = 2000
$ws.Range("L10").Value = 2000
$ws.Range("N10").Value = -2280
# Row 11 (item id G11=2481)
$ws.Range("H11").Value = 5089.5
$ws.Range("I11").Value = 186
$ws.Range("J11").Value = 19800
$ws.Range("K11").Value = 186
$ws.Range("L11").Value = 19800
$ws.Range("M11").Value = -46
$ws.Range("N11").Value = -20080
# Row 12 (item id G12=2392)
$ws.Range("H12").Value = 3600
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -332
# Row 37 (item id G37=2485)
$ws.Range("H37").Value = 5351
$ws.Range("I37").Value = 1651.6666
$ws.Range("J37").Value = 10900
$ws.Range("K37").Value = 1651.6666
$ws.Range("L37").Value = 10900
$ws.Range("M37").Value = -1514.6666
$ws.Range("N37").Value = -11174
# Row 80 (item id G80=13747)
$ws.Range("H80").Value = 39.47059
$ws.Range("I80").Value = 6
$ws.Range("J80").Value = 41.5625
$ws.Range("K80").Value = 6
$ws.Range("L80").Value = 41.5625
$ws.Range("M80").Value = 992
$ws.Range("N80").Value = -2037.5625
# Row 83 (item id G83=13747)
$ws.Range("H83").Value = 39.47059
$ws.Range("I83").Value = 6
$ws.Range("J83").Value = 41.5625
$ws.Range("K83").Value = 30
$ws.Range("L83").Value = 207.8125
$ws.Range("M83").Value = 4962
$ws.Range("N83").Value = -10191.8125
# Row 99 (item id G99=19943)
$ws.Range("H99").Value = 2063.3333
$ws.Range("I99").Value = 2300
$ws.Range("J99").Value = 1590
$ws.Range("K99").Value = 2300
$ws.Range("L99").Value = 1590
$ws.Range("M99").Value = -802
$ws.Range("N99").Value = -4586
# Row 112 (item id G112=25788)
$ws.Range("H112").Value = 35994.75
$ws.Range("J112").Value = 35994.75
$ws.Range("L112").Value = 35994.75
$ws.Range("N112").Value = -38948.75
# Row 134 (item id G134=43998)
$ws.Range("H134").Value = 2277.3157
$ws.Range("I134").Value = 1790.6923
$ws.Range("J134").Value = 3331.6667
$ws.Range("K134").Value = 5372.0769
$ws.Range("L134").Value = 9995.000100000001
$ws.Range("M134").Value = -2837.0769
$ws.Range("N134").Value = -15065.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 12 (item id G12=1604)
$ws.Range("H12").Value = 1774
$ws.Range("I12").Value = 1774
$ws.Range("K12").Value = 1774
$ws.Range("M12").Value = -1604
# Row 17 (item id G17=1823)
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("K17").Value = 8
$ws.Range("M17").Value = 166
# Row 31 (item id G31=44023)
$ws.Range("H31").Value = 3351362.2
$ws.Range("I31").Value = 7445454.5
$ws.Range("J31").Value = 1650.4849
$ws.Range("K31").Value = 7445454.5
$ws.Range("L31").Value = 1650.4849
$ws.Range("M31").Value = -7445159.5
$ws.Range("N31").Value = -2240.4849
# Row 34 (item id G34=44023)
$ws.Range("H34").Value = 3351362.2
$ws.Range("I34").Value = 7445454.5
$ws.Range("J34").Value = 1650.4849
$ws.Range("K34").Value = 7445454.5
$ws.Range("L34").Value = 1650.4849
$ws.Range("M34").Value = -7445252.5
$ws.Range("N34").Value = -2054.4849
# Row 107 (item id G107=27689)
$ws.Range("H107").Value = 1088.619
$ws.Range("I107").Value = 1190.7273
$ws.Range("J107").Value = 976.3
$ws.Range("K107").Value = 1190.7273
$ws.Range("L107").Value = 976.3
$ws.Range("M107").Value = 729.2727
$ws.Range("N107").Value = -4816.3

$ws = $wb.Worksheets.Item("CUL")
# Row 49 (item id G49=4719)
$ws.Range("H49").Value = 1200
$ws.Range("J49").Value = 1200
$ws.Range("L49").Value = 3600
$ws.Range("N49").Value = -3912
# Row 131 (item id G131=36060)
$ws.Range("H131").Value = 7937373
$ws.Range("J131").Value = 9009956
$ws.Range("L131").Value = 27029868
$ws.Range("N131").Value = -27039948

$ws = $wb.Worksheets.Item("GSM")
# Row 17 (item id G17=2445)
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2336

$ws = $wb.Worksheets.Item("LTW")
# Row 3 (item id G3=3537)
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 15 (item id G15=3537)
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 93 (item id G93=19993)
$ws.Range("H93").Value = 1129.0625
$ws.Range("I93").Value = 1224.1818
$ws.Range("J93").Value = 919.8
$ws.Range("K93").Value = 1224.1818
$ws.Range("L93").Value = 919.8
$ws.Range("M93").Value = 23.81819999999993
$ws.Range("N93").Value = -3415.8
# Row 136 (item id G136=44060)
$ws.Range("H136").Value = 5854.7407
$ws.Range("I136").Value = 11375.083
$ws.Range("J136").Value = 1438.4667
$ws.Range("K136").Value = 34125.249
$ws.Range("L136").Value = 4315.4001
$ws.Range("M136").Value = -31575.249
$ws.Range("N136").Value = -9415.400099999999

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (item id G100=19981)
$ws.Range("H100").Value = 7117.625
$ws.Range("I100").Value = 7705.857
$ws.Range("K100").Value = 15411.714
$ws.Range("M100").Value = -14870.714
